# Update the price list date and the six "Regaton de GOMA NEGRA" prices
# on the active sheet ("Hoja1"), matching the committed spreadsheet edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date (A1) moves forward one month (serial date 45406 -> 45436).
$ws.Range("A1").Value = 45436

# Updated unit prices in column D for rows 28-33.
$ws.Range("D28").Value = 8658
$ws.Range("D29").Value = 9139
$ws.Range("D30").Value = 10918
$ws.Range("D31").Value = 11980
$ws.Range("D32").Value = 12450
$ws.Range("D33").Value = 13280
